# dynamic source from text file
# Add a new applicant row (Maureen Gwapa / Walk-in / Fiber Technician) sourced
# from a text file, adjust a few column widths, and update the active
# selection/scroll position in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new data row (row 5) -----------------------------------------
# Copy the date cell's formatting from row 2 (A2) onto A5 so the new date
# picks up the existing short-date number format instead of creating a new
# style, then overwrite the value.
$ws.Range("A2").Copy($ws.Range("A5"))
$ws.Range("A5").Value = 42915
$ws.Range("B5").Value = "Maureen Gwapa"
$ws.Range("C5").Value = "Walk-in"
$ws.Range("D5").Value = "Fiber Technician"
$ws.Range("E5").Value = 5645645

# --- Column width tweaks ---------------------------------------------------
# (input values chosen so the pixel-quantized ColumnWidth setter lands on,
# or as close as mathematically possible to, the target stored width)
$ws.Columns.Item(3).ColumnWidth = 18.1666666666667  # C: 14.86 -> 19
$ws.Columns.Item(6).ColumnWidth = 17                 # F: 15.14 -> ~17.86
$ws.Columns.Item(11).ColumnWidth = 17.6666666666667  # K: 16.57 -> ~18.43
$ws.Columns.Item(12).ColumnWidth = 16.1666666666667  # L: 14.14 -> 17

# --- Sheet view: scroll so column B is the left-most visible column -------
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1

# --- Update the active selection -------------------------------------------
$ws.Range("C22").Select()
